$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update last-updated timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 15:36"

# Update country name order (A column) for rows whose rank changed after refresh/resort
# Row 4
$ws.Cells.Item(4, 2).Value = 2357323
$ws.Cells.Item(4, 3).Value = 666
$ws.Cells.Item(4, 4).Value = 980367
$ws.Cells.Item(4, 5).Value = 1254697
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 12
$ws.Cells.Item(4, 8).Value = 122259

# Row 7
$ws.Cells.Item(7, 2).Value = 427046
$ws.Cells.Item(7, 3).Value = 136
$ws.Cells.Item(7, 4).Value = 237929
$ws.Cells.Item(7, 5).Value = 175400
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 14
$ws.Cells.Item(7, 8).Value = 13717

# Row 14
$ws.Cells.Item(14, 2).Value = 191689
$ws.Cells.Item(14, 3).Value = 114
$ws.Cells.Item(14, 4).Value = 174900
$ws.Cells.Item(14, 5).Value = 7826
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 1
$ws.Cells.Item(14, 8).Value = 8963

# Row 18
$ws.Cells.Item(18, 1).Value = "Arabia Saudita"
$ws.Cells.Item(18, 2).Value = 161005
$ws.Cells.Item(18, 3).Value = 3393
$ws.Cells.Item(18, 4).Value = 105175
$ws.Cells.Item(18, 5).Value = 54523
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 40
$ws.Cells.Item(18, 8).Value = 1307

# Row 19
$ws.Cells.Item(19, 1).Value = "Francia"
$ws.Cells.Item(19, 2).Value = 160377
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 74372
$ws.Cells.Item(19, 5).Value = 56365
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 29640

# Row 23
$ws.Cells.Item(23, 2).Value = 88403
$ws.Cells.Item(23, 3).Value = 1034
$ws.Cells.Item(23, 4).Value = 69956
$ws.Cells.Item(23, 5).Value = 18348
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 1
$ws.Cells.Item(23, 8).Value = 99

# Row 34
$ws.Cells.Item(34, 2).Value = 42785
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(34, 4).Value = 13153
$ws.Cells.Item(34, 5).Value = 28616
$ws.Cells.Item(34, 6).Value = 0
$ws.Cells.Item(34, 7).Value = 5
$ws.Cells.Item(34, 8).Value = 1016

# Row 37
$ws.Cells.Item(37, 2).Value = 39392
$ws.Cells.Item(37, 3).Value = 259
$ws.Cells.Item(37, 4).Value = 25548
$ws.Cells.Item(37, 5).Value = 12310
$ws.Cells.Item(37, 6).Value = 0
$ws.Cells.Item(37, 7).Value = 4
$ws.Cells.Item(37, 8).Value = 1534

# Row 60
$ws.Cells.Item(60, 2).Value = 12990
$ws.Cells.Item(60, 3).Value = 96
$ws.Cells.Item(60, 4).Value = 11997
$ws.Cells.Item(60, 5).Value = 731
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 1
$ws.Cells.Item(60, 8).Value = 262

# Row 81
$ws.Cells.Item(81, 2).Value = 5196
$ws.Cells.Item(81, 3).Value = 90
$ws.Cells.Item(81, 4).Value = 1974
$ws.Cells.Item(81, 5).Value = 2975
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 7).Value = 9
$ws.Cells.Item(81, 8).Value = 247

# Row 85
$ws.Cells.Item(85, 1).Value = "Etiopia"
$ws.Cells.Item(85, 2).Value = 4663
$ws.Cells.Item(85, 3).Value = 131
$ws.Cells.Item(85, 4).Value = 1297
$ws.Cells.Item(85, 5).Value = 3291
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 75

# Row 86
$ws.Cells.Item(86, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(86, 2).Value = 4599
$ws.Cells.Item(86, 3).Value = 17
$ws.Cells.Item(86, 4).Value = 3952
$ws.Cells.Item(86, 5).Value = 599
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 3
$ws.Cells.Item(86, 8).Value = 48

# Row 92
$ws.Cells.Item(92, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(92, 2).Value = 3525
$ws.Cells.Item(92, 3).Value = 94
$ws.Cells.Item(92, 4).Value = 2270
$ws.Cells.Item(92, 5).Value = 1084
$ws.Cells.Item(92, 6).Value = 0
$ws.Cells.Item(92, 7).Value = 2
$ws.Cells.Item(92, 8).Value = 171

# Row 93
$ws.Cells.Item(93, 1).Value = "Kirguistan"
$ws.Cells.Item(93, 2).Value = 3356
$ws.Cells.Item(93, 3).Value = 205
$ws.Cells.Item(93, 4).Value = 2021
$ws.Cells.Item(93, 5).Value = 1295
$ws.Cells.Item(93, 6).Value = 0
$ws.Cells.Item(93, 7).Value = 3
$ws.Cells.Item(93, 8).Value = 40

# Row 100
$ws.Cells.Item(100, 2).Value = 2434
$ws.Cells.Item(100, 3).Value = 30
$ws.Cells.Item(100, 4).Value = 2188
$ws.Cells.Item(100, 5).Value = 215
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 31

# Row 107
$ws.Cells.Item(107, 1).Value = "Mali"
$ws.Cells.Item(107, 2).Value = 1961
$ws.Cells.Item(107, 3).Value = 28
$ws.Cells.Item(107, 4).Value = 1266
$ws.Cells.Item(107, 5).Value = 584
$ws.Cells.Item(107, 6).Value = 0
$ws.Cells.Item(107, 7).Value = 2
$ws.Cells.Item(107, 8).Value = 111

# Row 108
$ws.Cells.Item(108, 1).Value = "Sri Lanka"
$ws.Cells.Item(108, 2).Value = 1950
$ws.Cells.Item(108, 3).Value = 0
$ws.Cells.Item(108, 4).Value = 1526
$ws.Cells.Item(108, 5).Value = 413
$ws.Cells.Item(108, 6).Value = 0
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 11

# Row 115
$ws.Cells.Item(115, 1).Value = "Libano"
$ws.Cells.Item(115, 2).Value = 1603
$ws.Cells.Item(115, 3).Value = 16
$ws.Cells.Item(115, 4).Value = 1077
$ws.Cells.Item(115, 5).Value = 494
$ws.Cells.Item(115, 6).Value = 0
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = 32

# Row 116
$ws.Cells.Item(116, 1).Value = "Eslovaquia"
$ws.Cells.Item(116, 2).Value = 1588
$ws.Cells.Item(116, 3).Value = 1
$ws.Cells.Item(116, 4).Value = 1447
$ws.Cells.Item(116, 5).Value = 113
$ws.Cells.Item(116, 6).Value = 0
$ws.Cells.Item(116, 7).Value = 0
$ws.Cells.Item(116, 8).Value = 28

# Row 129
$ws.Cells.Item(129, 1).Value = "Estado de Palestina"
$ws.Cells.Item(129, 2).Value = 975
$ws.Cells.Item(129, 3).Value = 142
$ws.Cells.Item(129, 4).Value = 439
$ws.Cells.Item(129, 5).Value = 533
$ws.Cells.Item(129, 6).Value = 0
$ws.Cells.Item(129, 7).Value = 0
$ws.Cells.Item(129, 8).Value = 3

# Row 130
$ws.Cells.Item(130, 1).Value = "Yemen"
$ws.Cells.Item(130, 2).Value = 941
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 347
$ws.Cells.Item(130, 5).Value = 338
$ws.Cells.Item(130, 6).Value = 0
$ws.Cells.Item(130, 7).Value = 0
$ws.Cells.Item(130, 8).Value = 256

# Row 155
$ws.Cells.Item(155, 2).Value = 446
$ws.Cells.Item(155, 3).Value = 0
$ws.Cells.Item(155, 4).Value = 435
$ws.Cells.Item(155, 5).Value = 4
$ws.Cells.Item(155, 6).Value = 0
$ws.Cells.Item(155, 7).Value = 0
$ws.Cells.Item(155, 8).Value = 7

# Row 183
$ws.Cells.Item(183, 2).Value = 83
$ws.Cells.Item(183, 3).Value = 1
$ws.Cells.Item(183, 4).Value = 69
$ws.Cells.Item(183, 5).Value = 13
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 1

# Row 202
$ws.Cells.Item(202, 1).Value = "Dominica"

# Row 203
$ws.Cells.Item(203, 1).Value = "Fiyi"

# Row 211
$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 1

# Row 212
$ws.Cells.Item(212, 1).Value = "Seychelles"
$ws.Cells.Item(212, 2).Value = 11
$ws.Cells.Item(212, 3).Value = 0
$ws.Cells.Item(212, 4).Value = 11
$ws.Cells.Item(212, 5).Value = 0
$ws.Cells.Item(212, 6).Value = 0
$ws.Cells.Item(212, 7).Value = 0
$ws.Cells.Item(212, 8).Value = 0

# Row 214
$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 2).Value = 8
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1

# Row 215
$ws.Cells.Item(215, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(215, 2).Value = 8
$ws.Cells.Item(215, 3).Value = 0
$ws.Cells.Item(215, 4).Value = 8
$ws.Cells.Item(215, 5).Value = 0
$ws.Cells.Item(215, 6).Value = 0
$ws.Cells.Item(215, 7).Value = 0
$ws.Cells.Item(215, 8).Value = 0
